# Penalty Reward System changes (unfinished) - updates to po_data workbook
$wb = $excel.ActiveWorkbook

# --- "Weekly Quantity" sheet ---
# Remove the two weekly rows that correspond to 2023-09-17 and 2023-09-24
# (old rows 13 and 14), shifting all later rows up by two and shrinking
# the used range from A1:B21 to A1:B19.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("13:14").Delete()

# Update the requested quantity for the week of 2023-10-01 (row 11)
$wsWeekly.Range("B11").Value = 94

# --- "Monthly Trend" sheet ---
# Update the requested quantity for the month of 2023-09 (row 5)
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B5").Value = 178
